$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(515).Insert()

$ws.Cells.Item(515, 1).Value = 4
$ws.Cells.Item(515, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(515, 3).Value = "Los Lagos"
$ws.Cells.Item(515, 4).Value = 45267
$ws.Cells.Item(515, 5).Value = 10
$ws.Cells.Item(515, 6).Value = 100112043
$ws.Cells.Item(515, 7).Value = "Pepino ensalada"
$ws.Cells.Item(515, 8).Value = "Sin especificar"
$ws.Cells.Item(515, 9).Value = "Primera"
$ws.Cells.Item(515, 10).Value = 200
$ws.Cells.Item(515, 11).Value = 21000
$ws.Cells.Item(515, 12).Value = 21000
$ws.Cells.Item(515, 13).Value = 21000
$ws.Cells.Item(515, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(515, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(515, 16).Value = 350
$ws.Cells.Item(515, 17).Value = 60
$ws.Cells.Item(515, 18).Value = "Hortaliza"
